$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2022-09-05 07:12:02"
$newTimestamp = "2022-09-05 21:01:12"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 15).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 64 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 15)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value2 = $newTimestamp
    }
}
